$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 329, shifting existing rows 329:462 down to 330:463.
$ws.Rows("329:329").Insert()

# Populate the newly inserted row 329 with the new data point.
$ws.Range("A329").Value = 5
$ws.Range("B329").Value = "Macroferia Regional de Talca"
$ws.Range("C329").Value = "Maule"
$ws.Range("D329").Value = 45119
$ws.Range("E329").Value = 7
$ws.Range("F329").Value = 100112009
$ws.Range("G329").Value = "Acelga"
$ws.Range("H329").Value = "Sin especificar"
$ws.Range("I329").Value = "Primera"
$ws.Range("J329").Value = 500
$ws.Range("K329").Value = 1800
$ws.Range("L329").Value = 1800
$ws.Range("M329").Value = 1800
$ws.Range("N329").Value = "$/docena de atados (4 kilos)"
$ws.Range("O329").Value = "Región del Maule"
$ws.Range("P329").Value = 450
$ws.Range("Q329").Value = 4
$ws.Range("R329").Value = "Hortaliza"
